# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column N on the "Repayment Schedule"
#   sheet, shifting the old N/O/P ("Late" / "Heading" / "Outstanding")
#   columns one place to the right.
# - Make "Repayment Schedule" the active sheet/tab (was "Transactions"),
#   with a fresh selection.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column at N; everything from N onward shifts right.
$wsSchedule.Columns("N").Insert()

# The newly inserted column picks up the width of its left neighbour
# (column M) the way Excel's default insert behaviour does.
$wsSchedule.Columns("N").ColumnWidth = $wsSchedule.Columns("M").ColumnWidth

# Switch the active sheet to "Repayment Schedule" and move the selection,
# which also clears the previous selection/active-tab on "Transactions".
$wsSchedule.Select()
$wsSchedule.Range("R11").Select()
